$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '26.380.41'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +0.28%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.613.95'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  -0.02%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '213.20'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('E6').Value = '  -0.03%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.487'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.03%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.0616'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.55%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '18.54'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +2.16%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0814'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -1.10%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.838.50'
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.622.46'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('E14').Value = '  +0.25%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.516'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.46%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '26.358.29'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.35%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '62.02'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('E19').Value = '  -0.09%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '203.28'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.82%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '4.30'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.69%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '9.35'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.20%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '6.03'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('E24').Value = '  +6.75%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '144.47'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +1.47%  '
$ws.Range('E26').Value = '  +0.04%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.121'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -2.40%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '15.24'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.42%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '6.59'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.68%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0492'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +4.22%  '
$ws.Range('E31').Value = '  +0.31%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.21'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +1.87%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '2.95'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -1.72%  '
$ws.Range('E34').Value = '  +3.03%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.49'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.64%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.162.23'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +5.13%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.0166'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +1.46%  '
$ws.Range('E38').Value = '  -0.05%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.794'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.96%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.504'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.32'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.38%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.787'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +0.63%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '5.24'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.66%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.751.23'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('E45').Value = '  -1.55%  '
$ws.Range('E46').Value = '  -0.96%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '54.43'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.30%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0508'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0956'
$ws.Range('E50').Value = '  -11.16%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.03%  '
